$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 358
$ws.Range("F3").Value = 1279
$ws.Range("F5").Value = 2064
$ws.Range("F6").Value = 118
$ws.Range("F7").Value = 836
$ws.Range("F8").Value = 19
$ws.Range("F10").Value = 141
$ws.Range("F11").Value = 1065
$ws.Range("F12").Value = 804
$ws.Range("F14").Value = 664
$ws.Range("F15").Value = 1291
$ws.Range("F16").Value = 1023
$ws.Range("F18").Value = 744
$ws.Range("F19").Value = 721
$ws.Range("F20").Value = 78
$ws.Range("F21").Value = 578
$ws.Range("F22").Value = 96
$ws.Range("F23").Value = 634
$ws.Range("F24").Value = 1215
$ws.Range("F25").Value = 135
$ws.Range("F26").Value = 424
$ws.Range("F28").Value = 5103
$ws.Range("F29").Value = 238
$ws.Range("F31").Value = 2415
$ws.Range("F32").Value = 5795
$ws.Range("F33").Value = 121
$ws.Range("F34").Value = 966
$ws.Range("F35").Value = 584
$ws.Range("F36").Value = 58
$ws.Range("F38").Value = 1038
$ws.Range("F41").Value = 661

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 2083
$ws.Range("F8").Value = 111
$ws.Range("F9").Value = 463
$ws.Range("F13").Value = 122
$ws.Range("F15").Value = 654
$ws.Range("F16").Value = 654
$ws.Range("F23").Value = 17
$ws.Range("F28").Value = 1709
$ws.Range("F29").Value = 522
$ws.Range("F38").Value = 38
$ws.Range("F40").Value = 93
$ws.Range("F42").Value = 480
$ws.Range("F45").Value = 2
$ws.Range("F47").Value = 93
$ws.Range("F49").Value = 8

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 655
$ws.Range("F6").Value = 743
$ws.Range("F7").Value = 356
$ws.Range("F8").Value = 208

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 358
$ws.Range("F3").Value = 655
$ws.Range("F5").Value = 1280
$ws.Range("F7").Value = 356
$ws.Range("F8").Value = 208
$ws.Range("F9").Value = 208
$ws.Range("F10").Value = 2083
$ws.Range("F11").Value = 2064
$ws.Range("F13").Value = 836
$ws.Range("F14").Value = 111
$ws.Range("F15").Value = 19
$ws.Range("F17").Value = 141
$ws.Range("F18").Value = 1065
$ws.Range("F19").Value = 804
$ws.Range("F21").Value = 463
$ws.Range("F22").Value = 664
$ws.Range("F23").Value = 1291
$ws.Range("F25").Value = 744
$ws.Range("F26").Value = 122
$ws.Range("F27").Value = 721
$ws.Range("F28").Value = 78
$ws.Range("F29").Value = 578
$ws.Range("F30").Value = 654
$ws.Range("F31").Value = 634
$ws.Range("F32").Value = 1215
$ws.Range("F33").Value = 135
$ws.Range("F35").Value = 424
$ws.Range("F36").Value = 5103
$ws.Range("F37").Value = 238
$ws.Range("F39").Value = 2415
$ws.Range("F40").Value = 5795
$ws.Range("F41").Value = 966
$ws.Range("F42").Value = 1709
$ws.Range("F43").Value = 584
$ws.Range("F44").Value = 58
$ws.Range("F45").Value = 1038
$ws.Range("F46").Value = 661
$ws.Range("F49").Value = 38
$ws.Range("F51").Value = 480
